$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "rollback bimy to before migration": the K006 "placement count" KPI rows
# (15-18) get their display name (KPI Level 1/2 Name columns) restored to
# the pre-migration capitalised "Placement Count", while the underlying
# "Value" column keeps the lower-case "placement count" string.
$ws.Range("A15:B18").Value = "Placement Count"

# Restore the cursor/selection to where the author left it.
$ws.Range("E27").Select()
